# TC03_CDS_Filter_Study-Molecular Char Init.xlsx
# Update the "SamplesTab" Neo4j query (row 3, column B) on the startup sheet
# to return the sample's own tumor status directly instead of the collected
# list, and refresh the row height / selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTumorQuery = @"
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["Molecular Characterization Initiative"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as ``Sample ID``,
 coalesce(p.participant_id,'') as ``Participant ID``,
 coalesce(s.study_name, '') as ``Study Name``,
 coalesce(s.phs_accession,'') as ``Accession``,
 coalesce(samp.sample_tumor_status,'') as ``Tumor``,
coalesce(samp.sample_type,'') as ``Analyte Type``
  ORDER By samp.sample_id LIMIT 100
"@

# Strip the trailing newline the here-string adds after the last line.
$newTumorQuery = $newTumorQuery.TrimEnd("`r", "`n")

$ws.Range("B3").Value = $newTumorQuery

# The longer replacement text needs a taller row.
$ws.Rows.Item(3).RowHeight = 218.25

# Match the author's final active selection.
$ws.Range("B11").Select()

$wb.Save()
